$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 0.8106692983931111
$ws.Range("R2").Value = 7.296023685538
$ws.Range("S2").Value = 0.3324646254321348
$ws.Range("T2").Value = 0.3324646254321347

$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("S3").Value = 0.08209009580482188
$ws.Range("T3").Value = 0.08209009580482188

$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 0.7934299436544445
$ws.Range("R4").Value = 7.140869492889999
$ws.Range("S4").Value = 0.3253945715553649
$ws.Range("T4").Value = 0.3253945715553649

$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 0.3112379131508889
$ws.Range("R5").Value = 2.801141218358
$ws.Range("S5").Value = 0.1276421796422985
$ws.Range("T5").Value = 0.1276421796422985

$ws.Range("G6").Value = 0.05862833333333334
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.110264333333333
$ws.Range("N6").Value = 6.330793
$ws.Range("O6").Value = 0.3832041185227171
$ws.Range("P6").Value = 0.3832041185227171
$ws.Range("Q6").Value = 0.1237212807561111
$ws.Range("R6").Value = 1.113491526805
$ws.Range("S6").Value = 0.0507394930905823
$ws.Range("T6").Value = 0.0507394930905823

$ws.Range("G7").Value = 0.05862833333333334
$ws.Range("O7").Value = 0.0946183755984393
$ws.Range("P7").Value = 0.0946183755984393
$ws.Range("S7").Value = 0.01252827979361742
$ws.Range("T7").Value = 0.01252827979361742

$ws.Range("G8").Value = 0.05862833333333334
$ws.Range("M8").Value = 2.065388333333333
$ws.Range("N8").Value = 6.196165
$ws.Range("O8").Value = 0.3750550597762889
$ws.Range("P8").Value = 0.3750550597762889
$ws.Range("Q8").Value = 0.1210902756694445
$ws.Range("R8").Value = 1.089812481025
$ws.Range("S8").Value = 0.04966048822092397
$ws.Range("T8").Value = 0.04966048822092397

$ws.Range("G9").Value = 0.05862833333333334
$ws.Range("M9").Value = 0.8101876666666666
$ws.Range("N9").Value = 2.430563
$ws.Range("O9").Value = 0.1471224461025547
$ws.Range("P9").Value = 0.1471224461025547
$ws.Range("Q9").Value = 0.04749995258388889
$ws.Range("R9").Value = 0.427499573255
$ws.Range("S9").Value = 0.01948026646025624
$ws.Range("T9").Value = 0.01948026646025624
